$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Hydro"

$ws.Range("C13").Value = 3548.6364102564103
$ws.Range("C14").Value = 1310

$ws.Range("C5").Select()
